$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.374.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.885.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08084"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3131"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08357"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7214"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.299"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008467"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.379.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.130.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.824"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1590"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.081"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.507"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.431"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.340"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05377"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.955"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.182"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7508"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01882"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.283.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.745"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.572"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "110.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8927"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000129"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.027.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5213"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.498"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4374"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.71%  "
